$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (data rows 2-51) to Text format so numeric-looking strings are preserved exactly as text
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.430.17"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "2.239.16"
$ws.Range("E3").Value = "  -0.43%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "244.00"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").Value = "0.629"
$ws.Range("E6").Value = "  -0.34%  "
$ws.Range("D7").Value = "74.82"
$ws.Range("E7").Value = "  -2.31%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.609"
$ws.Range("E9").Value = "  -3.12%  "
$ws.Range("D10").Value = "43.16"
$ws.Range("E10").Value = "  -4.55%  "
$ws.Range("D11").Value = "0.0965"
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").Value = "7.03"
$ws.Range("E12").Value = "  -4.00%  "
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("D14").Value = "2.576.09"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").Value = "14.41"
$ws.Range("E15").Value = "  -2.38%  "
$ws.Range("D16").Value = "0.843"
$ws.Range("E16").Value = "  -2.54%  "
$ws.Range("D17").Value = "2.237.37"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "42.223.80"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("E19").Value = "  +4.77%  "
$ws.Range("D20").Value = "6.24"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("D21").Value = "73.21"
$ws.Range("E21").Value = "  +1.47%  "
$ws.Range("D22").Value = "11.24"
$ws.Range("E22").Value = "  +1.17%  "
$ws.Range("D23").Value = "231.58"
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").Value = "2.11"
$ws.Range("E24").Value = "  -6.81%  "
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").Value = "11.49"
$ws.Range("E26").Value = "  -3.98%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("D28").Value = "2.28"
$ws.Range("E28").Value = "  -1.78%  "
$ws.Range("E29").Value = "  -3.79%  "
$ws.Range("D30").Value = "167.18"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").Value = "20.66"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").Value = "5.72"
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").Value = "0.0805"
$ws.Range("E33").Value = "  -2.69%  "
$ws.Range("D34").Value = "30.37"
$ws.Range("E34").Value = "  -6.42%  "
$ws.Range("E35").Value = "  -0.35%  "
$ws.Range("E36").Value = "  -8.88%  "
$ws.Range("D37").Value = "4.37"
$ws.Range("E37").Value = "  -7.24%  "
$ws.Range("D38").Value = "0.0305"
$ws.Range("E38").Value = "  -4.13%  "
$ws.Range("D39").Value = "13.67"
$ws.Range("E39").Value = "  -5.12%  "
$ws.Range("D40").Value = "2.15"
$ws.Range("E40").Value = "  -2.21%  "
$ws.Range("D43").Value = "0.200"
$ws.Range("E43").Value = "  -1.62%  "
$ws.Range("D44").Value = "8.77"
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("D45").Value = "105.04"
$ws.Range("E45").Value = "  -3.00%  "
$ws.Range("E46").Value = "  -2.32%  "
$ws.Range("E47").Value = "  -1.60%  "
$ws.Range("E48").Value = "  -1.70%  "
$ws.Range("E49").Value = "  -1.56%  "
$ws.Range("D50").Value = "2.69"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("D51").Value = "2.447.11"
$ws.Range("E51").Value = "  -0.68%  "

# Row 41/42: MultiversX and THORChain swap positions with updated values
$ws.Range("B41").Value = "THORChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D41").Value = "5.72"
$ws.Range("E41").Value = "  -1.52%  "

$ws.Range("B42").Value = "MultiversX"
$ws.Range("C42").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D42").Value = "65.33"
$ws.Range("E42").Value = "  +1.86%  "

# Restore default (Normal) style on column D so no extra style id is introduced
$ws.Range("D2:D51").Style = "Normal"
